$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.453.97'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '1.871.46'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = '''243.77'
$ws.Range("E5").Value = '  +0.31%  '
$ws.Range("D6").Value = '''0.7058'
$ws.Range("E6").Value = '  -2.42%  '
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").Value = '''0.07950'
$ws.Range("E8").Value = '  -0.81%  '
$ws.Range("D9").Value = '''0.3146'
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").Value = '''24.53'
$ws.Range("E10").Value = '  -1.98%  '
$ws.Range("D11").Value = '''0.07814'
$ws.Range("E11").Value = '  -4.56%  '
$ws.Range("D12").Value = '1.893.28'
$ws.Range("E12").Value = '  +1.05%  '
$ws.Range("D13").Value = '''93.90'
$ws.Range("E13").Value = '  -0.97%  '
$ws.Range("D14").Value = '''5.181'
$ws.Range("E14").Value = '  -1.13%  '
$ws.Range("D15").Value = '''0.7037'
$ws.Range("E15").Value = '  -1.41%  '
$ws.Range("D16").Value = '''6.494'
$ws.Range("E16").Value = '  +0.94%  '
$ws.Range("D17").Value = '''0.000008588'
$ws.Range("E17").Value = '  +0.97%  '
$ws.Range("D18").Value = '29.492.85'
$ws.Range("E18").Value = '  +0.43%  '
$ws.Range("D19").Value = '''252.86'
$ws.Range("E19").Value = '  +3.38%  '
$ws.Range("D20").Value = '2.151.24'
$ws.Range("E20").Value = '  +0.74%  '
$ws.Range("E21").Value = '  -1.59%  '
$ws.Range("D22").Value = '''0.9995'
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").Value = '''7.643'
$ws.Range("E23").Value = '  -1.55%  '
$ws.Range("D24").Value = '''0.9999'
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("D25").Value = '''0.1549'
$ws.Range("E25").Value = '  -3.74%  '
$ws.Range("D26").Value = '''9.021'
$ws.Range("E26").Value = '  -0.36%  '
$ws.Range("D27").Value = '''161.39'
$ws.Range("E27").Value = '  -0.88%  '
$ws.Range("D28").Value = '''18.79'
$ws.Range("E28").Value = '  +1.38%  '
$ws.Range("D29").Value = '''1.497'
$ws.Range("E29").Value = '  -0.55%  '
$ws.Range("D30").Value = '''4.315'
$ws.Range("E30").Value = '  -2.16%  '
$ws.Range("D31").Value = '''4.271'
$ws.Range("E31").Value = '  -0.34%  '
$ws.Range("D32").Value = '''1.206'
$ws.Range("E32").Value = '  -2.02%  '
$ws.Range("D33").Value = '''0.05293'
$ws.Range("E33").Value = '  -1.35%  '
$ws.Range("D34").Value = '''1.901'
$ws.Range("E34").Value = '  -2.07%  '
$ws.Range("D35").Value = '''0.7619'
$ws.Range("E35").Value = '  -0.45%  '
$ws.Range("D36").Value = '''1.191'
$ws.Range("E36").Value = '  +0.90%  '
$ws.Range("D37").Value = '''2.702'
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("D38").Value = '''0.01883'
$ws.Range("E38").Value = '  +0.53%  '
$ws.Range("D39").Value = '1.279.64'
$ws.Range("E39").Value = '  +1.21%  '
$ws.Range("E40").Value = '  +0.49%  '
$ws.Range("D41").Value = '''0.9003'
$ws.Range("D42").Value = '''109.96'
$ws.Range("E42").Value = '  -3.29%  '
$ws.Range("D43").Value = '''6.010'
$ws.Range("E43").Value = '  -6.67%  '
$ws.Range("D44").Value = '''71.05'
$ws.Range("E44").Value = '  -4.69%  '
$ws.Range("D45").Value = '''0.9997'
$ws.Range("E45").Value = '  -0.22%  '
$ws.Range("D46").Value = '2.047.79'
$ws.Range("E46").Value = '  +1.20%  '
$ws.Range("D47").Value = '''0.00000000127'
$ws.Range("E47").Value = '  -3.03%  '
$ws.Range("D48").Value = '''9.662'
$ws.Range("E48").Value = '  +1.66%  '
$ws.Range("D49").Value = '''1.805'
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("D50").Value = '''0.5176'
$ws.Range("E50").Value = '  -0.49%  '
$ws.Range("D51").Value = '''0.4306'
$ws.Range("E51").Value = '  -0.95%  '
